$wb = $excel.ActiveWorkbook

function Get-HyperlinkAtAddress($ws, $addrWanted) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addrWanted) {
            return $hl
        }
    }
    return $null
}

$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f2c776f39f6a9a77e8607472ae6b83c922c1ba4/e2e/54bee925-7da0-4db9-9d2a-c051257bf26a.md"
$latestUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af4db6209bc3418bf3192ef7c4021cc805733e74/e2e/54bee925-7da0-4db9-9d2a-c051257bf26a.md"
$errorDetail = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."
$mdDisplay = "54bee925-7da0-4db9-9d2a-c051257bf26a.md"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen the Error Detail column (P / column 16) to fit the new message.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Preserve the existing hyperlinks on A7/A8 so they can be re-inserted
    # after the new I6 hyperlink (keeps relationship id ordering correct).
    $a7hl = Get-HyperlinkAtAddress $ws "`$A`$7"
    $a7Address = $a7hl.Address
    $a7Display = $a7hl.TextToDisplay
    $a7hl.Delete()

    $a8hl = Get-HyperlinkAtAddress $ws "`$A`$8"
    $a8Address = $a8hl.Address
    $a8Display = $a8hl.TextToDisplay
    $a8hl.Delete()

    # Row 6 corresponds to 54bee925-7da0-4db9-9d2a-c051257bf26a (Ready for handoff).
    # A handback file has now been generated for it, so populate the
    # "Latest Target File", "Latest Handback File", "Latest Handback DateTime"
    # and "Error Detail" columns (I, J, K, P).
    $ws.Range("I6").Value = $mdDisplay
    $ws.Range("J6").Value = $ws.Range("G6").Value
    if ($sheetName -eq "zh-cn") {
        $ws.Range("K6").Value = "2016-08-30 20:48:44"
    } else {
        $ws.Range("K6").Value = "2016-08-30 20:48:51"
    }
    $ws.Range("P6").Value = $errorDetail

    # Re-create the hyperlinks in the expected order: ..., A6, I6, A7, A8
    $ws.Hyperlinks.Add($ws.Range("I6"), $currentUrl, "", "", $mdDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A7"), $a7Address, "", "", $a7Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A8"), $a8Address, "", "", $a8Display) | Out-Null
}
